# problem 2 fix bug
# Relabel the cluster ids in column A (rows 2..146) according to the
# permutation 0->3, 1->2, 2->0, 3->4, 4->1. Row 1 (header/first value,
# which carries the bold/bordered style) is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$map = @{ 0 = 3; 1 = 2; 2 = 0; 3 = 4; 4 = 1 }

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 146 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $old = [int]$cell.Value2
    $cell.Value2 = $map[$old]
}
